# Add new columns I ("I0") and J ("IF") with their data, mirroring the
# existing header formatting used by columns B:H (bold, bordered, centered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) from the existing header cell H1 so the new
# headers look consistent with the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for rows 2-19
$values = @(
    @(2, 6),
    @(2, 7),
    @(2, 6),
    @(3, 7),
    @(3, 6),
    @(1, 3),
    @(1, 6),
    @(2, 7),
    @(2, 6),
    @(4, 8),
    @(1, 4),
    @(6, 6),
    @(1, 3),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(6, 6)
)

$row = 2
foreach ($pair in $values) {
    $ws.Range("I$row").Value = $pair[0]
    $ws.Range("J$row").Value = $pair[1]
    $row++
}
